$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new documentation sections ("ORDER DISPLAY TABLE FOR AGENT
# (SHORT)" and "ORDER DISPLAY TABLE FOR AGENT (DETAIL)") right after the
# "AGENT DETAIL TABLE" block (which ends at row 24) and before the existing
# "TRACKING TABLE" block (which starts at old row 26). This pushes every
# row from the old row 26 onward down by 6 rows.
# ---------------------------------------------------------------------------
$ws.Range("A26:A31").EntireRow.Insert()

# --- Section: ORDER DISPLAY TABLE FOR AGENT(SHORT) -> rows 26 (header) & 27 (fields) ---
$ws.Range("B17:N17").Copy()
$ws.Range("B26:N26").PasteSpecial(-4122)
$ws.Rows.Item(26).RowHeight = 15.75
$ws.Range("B26").Value = "ORDER DISPLAY TABLE FOR AGENT(SHORT)"
$ws.Range("B26:N26").Merge()

$ws.Range("B27").Value = "order_id"
$ws.Range("C27").Value = "order_tracking_id"
$ws.Range("D27").Value = "order_delivery_state"
$ws.Range("E27").Value = "order_delivery_city"
$ws.Range("F27").Value = "order_place_date"

$ws.Rows.Item(28).RowHeight = 15.75

# --- Section: ORDER DISPLAY TABLE FOR AGENT (DETAIL) -> rows 29 (header) & 30 (fields) ---
$ws.Range("B17:N17").Copy()
$ws.Range("B29:N29").PasteSpecial(-4122)
$ws.Rows.Item(29).RowHeight = 15.75
$ws.Range("B29").Value = "ORDER DISPLAY TABLE FOR AGENT (DETAIL)"
$ws.Range("B29:N29").Merge()

$ws.Range("B30").Value = "order_id"
$ws.Range("C30").Value = "order_tracking_id "
$ws.Range("D30").Value = "order_delivery_name"
$ws.Range("E30").Value = "order_delivery_contact"
$ws.Range("F30").Value = "order_place_date"
$ws.Range("G30").Value = "order_delivery_address"
$ws.Range("H30").Value = "order_delivery_city"
$ws.Range("I30").Value = "order_delivery_state"

$ws.Rows.Item(31).RowHeight = 15.75

# Restore the on-screen selection to the new field row, matching where the
# author was working when the sections were added.
$ws.Range("E27").Select()

Write-Output "done"
